# cetak telaah penambahan tanda tangan
# Update the "telaah" header block from a "rencana pengadaan" report
# to a "rencana pemeliharaan" report, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HASIL PENELAAHAN RENCANA PEMELIHARAAN BARANG MILIK DAERAH"
$ws.Range("A2").Value = "(RENCANA PEMELIHARAAN)"
$ws.Range("A3").Value = "PEMERINTAH KABUPATEN BANJARNEGARA"
$ws.Range("A4").Value = ""

$ws.Range("H5").Select()
